$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 23533.334
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 32800
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 32800
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -34048

$ws.Range("H65").Value = 23533.334
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 32800
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 164000
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -170240

$ws.Range("H86").Value = 64015676
$ws.Range("I86").Value = 80018740
$ws.Range("J86").Value = 3440.8
$ws.Range("K86").Value = 80018740
$ws.Range("L86").Value = 3440.8
$ws.Range("M86").Value = -80017617
$ws.Range("N86").Value = -5686.8

$ws.Range("H89").Value = 64015676
$ws.Range("I89").Value = 80018740
$ws.Range("J89").Value = 3440.8
$ws.Range("K89").Value = 400093700
$ws.Range("L89").Value = 17204
$ws.Range("M89").Value = -400088084
$ws.Range("N89").Value = -28436

$ws.Range("H98").Value = 1024.6666
$ws.Range("I98").Value = 938.7826
$ws.Range("K98").Value = 938.7826
$ws.Range("M98").Value = 559.2174

$ws.Range("H106").Value = 46158620
$ws.Range("I106").Value = 85720800
$ws.Range("J106").Value = 2750
$ws.Range("K106").Value = 85720800
$ws.Range("L106").Value = 2750
$ws.Range("M106").Value = -85720169
$ws.Range("N106").Value = -4012

$ws.Range("H107").Value = 1092.4783
$ws.Range("I107").Value = 1217.7368
$ws.Range("J107").Value = 497.5
$ws.Range("K107").Value = 1217.7368
$ws.Range("L107").Value = 497.5
$ws.Range("M107").Value = 702.2632000000001
$ws.Range("N107").Value = -4337.5

$ws.Range("H122").Value = 1024.6666
$ws.Range("I122").Value = 938.7826
$ws.Range("K122").Value = 2816.3478
$ws.Range("M122").Value = -366.3478

$ws.Range("H125").Value = 1031.1428
$ws.Range("J125").Value = 1031.1428
$ws.Range("L125").Value = 9280.2852
$ws.Range("N125").Value = -14200.2852

$ws.Range("H127").Value = 1508.4
$ws.Range("I127").Value = 500
$ws.Range("K127").Value = 1500
$ws.Range("M127").Value = 3460

$ws.Range("H132").Value = 3008.7058
$ws.Range("I132").Value = 3046.4482
$ws.Range("J132").Value = 2789.8
$ws.Range("K132").Value = 9139.3446
$ws.Range("L132").Value = 8369.400000000001
$ws.Range("M132").Value = -6609.3446
$ws.Range("N132").Value = -13429.4

$ws.Range("H141").Value = 5097.65
$ws.Range("I141").Value = 2489.0833
$ws.Range("J141").Value = 9010.5
$ws.Range("K141").Value = 7467.249899999999
$ws.Range("L141").Value = 27031.5
$ws.Range("M141").Value = -2287.249899999999
$ws.Range("N141").Value = -37391.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2292.1562
$ws.Range("I74").Value = 1458.0869
$ws.Range("J74").Value = 4423.6665
$ws.Range("K74").Value = 1458.0869
$ws.Range("L74").Value = 4423.6665
$ws.Range("M74").Value = -584.0869
$ws.Range("N74").Value = -6171.6665

$ws.Range("H77").Value = 2292.1562
$ws.Range("I77").Value = 1458.0869
$ws.Range("J77").Value = 4423.6665
$ws.Range("K77").Value = 7290.4345
$ws.Range("L77").Value = 22118.3325
$ws.Range("M77").Value = -2922.4345
$ws.Range("N77").Value = -30854.3325

$ws.Range("H122").Value = 251500
$ws.Range("I122").Value = 334666.66
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 1003999.98
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1001549.98
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4472.9683
$ws.Range("I31").Value = 1295.0975
$ws.Range("J31").Value = 6885.7964
$ws.Range("K31").Value = 1295.0975
$ws.Range("L31").Value = 6885.7964
$ws.Range("M31").Value = -1000.0975
$ws.Range("N31").Value = -7475.7964

$ws.Range("H34").Value = 4472.9683
$ws.Range("I34").Value = 1295.0975
$ws.Range("J34").Value = 6885.7964
$ws.Range("K34").Value = 1295.0975
$ws.Range("L34").Value = 6885.7964
$ws.Range("M34").Value = -1093.0975
$ws.Range("N34").Value = -7289.7964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4001980
$ws.Range("I4").Value = 14000530
$ws.Range("J4").Value = 2560
$ws.Range("K4").Value = 42001590
$ws.Range("L4").Value = 7680
$ws.Range("M4").Value = -42001478
$ws.Range("N4").Value = -7904

$ws.Range("H12").Value = 65.09524
$ws.Range("I12").Value = 35.18182
$ws.Range("K12").Value = 105.54546
$ws.Range("M12").Value = 67.45453999999999

$ws.Range("H126").Value = 4436
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4436
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 13308
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -23188

$ws.Range("H130").Value = 2007.3077
$ws.Range("I130").Value = 500
$ws.Range("J130").Value = 2281.3635
$ws.Range("K130").Value = 1500
$ws.Range("L130").Value = 6844.0905
$ws.Range("M130").Value = 3520
$ws.Range("N130").Value = -16884.0905

$ws.Range("H131").Value = 2893.2205
$ws.Range("I131").Value = 463.07693
$ws.Range("J131").Value = 3580
$ws.Range("K131").Value = 1389.23079
$ws.Range("L131").Value = 10740
$ws.Range("M131").Value = 3650.76921
$ws.Range("N131").Value = -20820

$ws.Range("H141").Value = 9325.950000000001
$ws.Range("I141").Value = 8156.273
$ws.Range("K141").Value = 24468.819
$ws.Range("M141").Value = -19288.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1318.9412
$ws.Range("I102").Value = 1334.8
$ws.Range("K102").Value = 1334.8
$ws.Range("M102").Value = 287.2

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4051
$ws.Range("I40").Value = 4085
$ws.Range("K40").Value = 4085
$ws.Range("M40").Value = -3949

$ws.Range("H61").Value = 4205.8125
$ws.Range("I61").Value = 3214.6667
$ws.Range("K61").Value = 3214.6667
$ws.Range("M61").Value = -3012.6667

$ws.Range("H88").Value = 38400
$ws.Range("J88").Value = 38000
$ws.Range("L88").Value = 38000
$ws.Range("N88").Value = -38856

$ws.Range("H91").Value = 38400
$ws.Range("J91").Value = 38000
$ws.Range("L91").Value = 38000
$ws.Range("N91").Value = -40964

$ws.Range("H100").Value = 66689.92999999999
$ws.Range("I100").Value = 71418.08
$ws.Range("K100").Value = 71418.08
$ws.Range("M100").Value = -70877.08

$ws.Range("H113").Value = 4205.8125
$ws.Range("I113").Value = 3214.6667
$ws.Range("K113").Value = 3214.6667
$ws.Range("M113").Value = -1044.6667

$ws.Range("H122").Value = 5990.4
$ws.Range("I122").Value = 4952
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 14856
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -12406
$ws.Range("N122").Value = -23650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1123.1052
$ws.Range("J113").Value = 556.7778
$ws.Range("L113").Value = 1670.3334
$ws.Range("N113").Value = -6010.3334

$ws.Range("H122").Value = 2127.6
$ws.Range("I122").Value = 2167.889
$ws.Range("J122").Value = 1765
$ws.Range("K122").Value = 6503.667
$ws.Range("L122").Value = 5295
$ws.Range("M122").Value = -4053.667
$ws.Range("N122").Value = -10195

Write-Host "Applied Anima_Profits scheduled-runner updates"
